# =====================================================================
# Edit script: upload new Part 2 and test macro workbooks
#  1. Add a new "Sheet3" worksheet after Sheet2, populated with book data
#  2. Fix the "Saleperson" typo -> "Salesperson" header on Sheet1
#  3. Update Sheet1 selection to B1
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---- Step 1: add Sheet3 after the last existing sheet -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

# Header row
$ws3.Cells.Item(1, 1).Value = "Title"
$ws3.Cells.Item(1, 2).Value = "Title_Remark"
$ws3.Cells.Item(1, 3).Value = "Author"
$ws3.Cells.Item(1, 4).Value = "Publisher"
$ws3.Cells.Item(1, 5).Value = "Pub_Year"
$ws3.Cells.Item(1, 6).Value = "ISBN"
$ws3.Cells.Item(1, 7).Value = "Binding"

# Book data (column B / Title_Remark intentionally left blank)
$bookTitles = @(
    "COLLECTIVE MEMORY WORK: A METHODOLOGY FOR LEARNING WITH AND FROM LIVED EXPERIENCE",
    "COMPASSIONATE MIGRATION AND REGIONAL POLICY IN THE AMERICAS",
    "GENDERING NATIONALISM: INTERSECTIONS OF NATION, GENDER AND SEXUALITY",
    "GRIDDED WORLDS: AN URBAN ANTHOLOGY",
    "HANDBOOK OF COMMUNITY MOVEMENTS AND LOCAL ORGANIZATIONS IN THE 21ST CENTURY",
    "IMAGE POLITICS OF CLIMATE CHANGE: VISUALIZATIONS, IMAGINATIONS, DOCUMENTATIONS",
    "IMMIGRATION AND METROPOLITAN REVITALIZATION IN THE UNITED ST.",
    "IMMIGRATION POLICY IN THE AGE OF PUNISHMENT: DETENTION, DEPORTATION, AND BORDER CONTROL.",
    "LAND RIGHTS, BIODIVERSITY CONSERVATION AND JUSTICE: RETHINKING PARKS AND PEOPLE",
    "MACHINE LEARNING TECHNIQUES FOR ONLINE SOCIAL NETWORKS",
    "MODERN AMERICAN EXTREMISM AND DOMESTIC TERRORISM: AN ENCYCLOPEDIA OF EXTREMISTS AND EXTREMIST GROUPS.",
    "ON REPLACEMENT: CULTURAL, SOCIAL AND PSYCHOLOGICAL REPRESENTATIONS",
    "ORGANIZING NETWORKS: AN ACTOR-NETWORK THEORY OF ORGANIZATIONS.",
    "POLITICS, POWER AND COMMUNITY DEVELOPMENT",
    "REVEALING TACIT KNOWLEDGE: EMBODIMENT AND EXPLICATION",
    "RURAL POVERTY IN THE UNITED STATES",
    "SHARED PROSPERITY IN AMERICA'S COMMUNITIES",
    "SOCIAL THEORIES OF URBAN VIOLENCE IN THE GLOBAL SOUTH: TOWARDS SAFE AND INCLUSIVE CITIES",
    "THICK SPACE.",
    "URBAN TRANSFORMATIONS IN THE U.S.A.: SPACES, COMMUNITIES, REPRESENTATIONS",
    "WHAT IS A SLAVE SOCIETY?: THE PRACTICE OF SLAVERY IN GLOBAL PERSPECTIVE",
    "WILDLIFE CRIME: FROM THEORY TO PRACTICE",
    "WORLDWIDE MOBILIZATIONS: CLASS STRUGGLES AND URBAN COMMONING"
)
$bookAuthors = @(
    "COREY W. JOHNSON",
    "STEVEN W. BENDER",
    "JON MULHOLLAND",
    "REUBEN ROSE-REDWOOD",
    "RAM A. CNAAN",
    "BIRGIT SCHNEIDER",
    "DOMENIC VITIELLO",
    "DAVID C. BROTHERTON",
    "SHARLENE MOLLETT",
    "TANSEL OZYER",
    "BARRY J. BALLECK",
    "JEAN OWEN",
    "ANDREA BELLINGER",
    "ROSIE MEADE",
    "FRANK ADLOFF",
    "ANN R. TICKAMYER",
    "SUSAN M. WACHTER",
    "JENNIFER ERIN SALAHUB",
    "DOROTHEE BRANTZ",
    "JULIA SATTLER",
    "NOEL LENSKI",
    "WILLIAM D. MORETO",
    "DON KALB"
)
$bookPublishers = @(
    "ROUTLEDGE",
    "PALGRAVE MACMILLAN",
    "PALGRAVE MACMILLAN",
    "SPRINGER",
    "SPRINGER",
    "TRANSCRIPT",
    "UNIV PENN PRESS",
    "COLUMBIA UNIVERSITY PRESS",
    "ROUTLEDGE",
    "SPRINGER",
    "ABC-CLIO",
    "PALGRAVE MACMILLAN",
    "TRANSCRIPT",
    "POLICY PRESS",
    "TRANSCRIPT VERLAG",
    "COLUMBIA UNIVERSITY PRESS",
    "UNIV OF PENNSYLVANIA PR",
    "ROUTLEDGE",
    "TRANSCRIPT-VERLAG",
    "TRANSCRIPT VERLAG",
    "CAMBRIDGE UNIV PRESS",
    "TEMPLE UNIVERSITY PRESS",
    "BERGHAHN BOOKS"
)
$bookYears = @(
    2018,
    2017,
    2018,
    2018,
    2018,
    2014,
    2017,
    2018,
    2018,
    2018,
    2018,
    2018,
    2016,
    2016,
    2015,
    2017,
    2016,
    2018,
    2014,
    2016,
    2018,
    2018,
    2018
)
$bookIsbns = @(
    9781315298696,
    9781137550743,
    9783319766997,
    9783319764900,
    9783319774169,
    9783839426104,
    9780812293951,
    9780231545891,
    9781315439464,
    9783319899329,
    9781440852756,
    9783319760117,
    9783839436165,
    9781447317388,
    9783839425169,
    9780231544719,
    9780812292404,
    9781351254700,
    9783839420430,
    9783839431115,
    9781108633208,
    9781439914731,
    9781785339073
)
$bookBindings = @(
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook",
    "eBook"
)

for ($i = 0; $i -lt $bookTitles.Length; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 1).Value = $bookTitles[$i]
    $ws3.Cells.Item($row, 3).Value = $bookAuthors[$i]
    $ws3.Cells.Item($row, 4).Value = $bookPublishers[$i]
    $ws3.Cells.Item($row, 5).Value = $bookYears[$i]
    $ws3.Cells.Item($row, 6).Value = $bookIsbns[$i]
    $ws3.Cells.Item($row, 7).Value = $bookBindings[$i]
}

# Column B (Title_Remark) display width, as in the source workbook
$ws3.Columns.Item(2).ColumnWidth = 16

# ---- Step 2: fix the Sheet1 header typo --------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B1").Value = "Salesperson"

# ---- Step 3: move the Sheet1 selection to B1 ---------------------------
$ws1.Activate() | Out-Null
$ws1.Range("B1").Select() | Out-Null

